$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two employee/character names between G60 and G91
# (underlying shared-string table reorder: "Eeth Koth" moves to sit
#  right after "Obi-Wan Kenobi" / "Palo Jemabie" moves to the end of that block)
$ws.Range("G60").Value = "Eeth Koth"
$ws.Range("G91").Value = "Palo Jemabie"

# Apply updated numeric cell values
$ws.Range("E2").Value = [double]"0.02144550121794863"
$ws.Range("E3").Value = [double]"2.455480850087047e-14"
$ws.Range("E4").Value = [double]"0.0585334811222988"
$ws.Range("E5").Value = [double]"0.006632132612788942"
$ws.Range("E8").Value = [double]"2.455480850087047e-14"
$ws.Range("E10").Value = [double]"2.455480850087047e-14"
$ws.Range("E11").Value = [double]"0.2618233098734618"
$ws.Range("E13").Value = [double]"2.455480850087047e-14"
$ws.Range("E14").Value = [double]"2.514412390489136e-11"
$ws.Range("E15").Value = [double]"0.01247717905899594"
$ws.Range("E18").Value = [double]"0.09717836577847853"
$ws.Range("E20").Value = [double]"0.04686379776670074"
$ws.Range("E21").Value = [double]"0.1482361462768889"
$ws.Range("E22").Value = [double]"2.687769338505282e-10"
$ws.Range("E23").Value = [double]"0.2982687018315968"
$ws.Range("B24").Value = [double]"0.01642789718572909"
$ws.Range("E27").Value = [double]"1.866165446066156e-12"
$ws.Range("E30").Value = [double]"0.07287551024603613"
$ws.Range("E32").Value = [double]"2.455480850087047e-14"
$ws.Range("E35").Value = [double]"2.455480850087047e-14"
$ws.Range("E37").Value = [double]"2.455480850087047e-14"
$ws.Range("E38").Value = [double]"0.09783422281736534"
$ws.Range("E39").Value = [double]"0.212089909069801"
$ws.Range("E40").Value = [double]"2.455480850087047e-14"
$ws.Range("E41").Value = [double]"0.003900711649066044"
$ws.Range("E43").Value = [double]"2.455480850087047e-14"
$ws.Range("E46").Value = [double]"2.455480850087047e-14"
$ws.Range("E48").Value = [double]"2.455480850087047e-14"
$ws.Range("E49").Value = [double]"2.455480850087047e-14"
$ws.Range("E50").Value = [double]"0.06303989561971357"
$ws.Range("E54").Value = [double]"2.455480850087047e-14"
$ws.Range("E55").Value = [double]"2.455480850087047e-14"
$ws.Range("E56").Value = [double]"2.455480850087047e-14"
$ws.Range("E57").Value = [double]"0.08262723622105826"
$ws.Range("B60").Value = [double]"0.0004826622335880557"
$ws.Range("C60").Value = [double]"16"
$ws.Range("D60").Value = [double]"0.1194029850746269"
$ws.Range("E60").Value = [double]"0.07376956536638371"
$ws.Range("F60").Value = [double]"9"
$ws.Range("H60").Value = [double]"7"
$ws.Range("E61").Value = [double]"2.455480850087047e-14"
$ws.Range("E63").Value = [double]"0.2232797699885382"
$ws.Range("E64").Value = [double]"0.0864203207321078"
$ws.Range("E65").Value = [double]"2.514412390489136e-11"
$ws.Range("E66").Value = [double]"2.455480850087047e-14"
$ws.Range("E67").Value = [double]"2.455480850087047e-14"
$ws.Range("E69").Value = [double]"2.455480850087047e-14"
$ws.Range("E70").Value = [double]"0.003900711649066044"
$ws.Range("E71").Value = [double]"3.769163104883618e-11"
$ws.Range("E72").Value = [double]"2.455480850087047e-14"
$ws.Range("E73").Value = [double]"3.771618585733705e-11"
$ws.Range("E75").Value = [double]"0.01713462530082498"
$ws.Range("E76").Value = [double]"2.455480850087047e-14"
$ws.Range("B77").Value = [double]"0.009426707872971615"
$ws.Range("E78").Value = [double]"0.06411431120940382"
$ws.Range("E79").Value = [double]"0.06474628483253868"
$ws.Range("E80").Value = [double]"2.455480850087047e-14"
$ws.Range("E81").Value = [double]"2.514412390489136e-11"
$ws.Range("E82").Value = [double]"0.06461640467758835"
$ws.Range("E85").Value = [double]"0.04530344405071713"
$ws.Range("E88").Value = [double]"2.455480850087047e-14"
$ws.Range("E89").Value = [double]"2.455480850087047e-14"
$ws.Range("E90").Value = [double]"2.455480850087047e-14"
$ws.Range("B91").Value = [double]"0"
$ws.Range("C91").Value = [double]"4"
$ws.Range("D91").Value = [double]"0.02985074626865672"
$ws.Range("E91").Value = [double]"0.01561186977399064"
$ws.Range("F91").Value = [double]"1"
$ws.Range("H91").Value = [double]"3"
$ws.Range("E92").Value = [double]"2.455480850087047e-14"
$ws.Range("E93").Value = [double]"2.455480850087047e-14"
$ws.Range("B95").Value = [double]"0.003996554517224881"
$ws.Range("E95").Value = [double]"0.1301688562996523"
$ws.Range("E96").Value = [double]"2.455480850087047e-14"
$ws.Range("E98").Value = [double]"2.455480850087047e-14"
$ws.Range("E102").Value = [double]"0.05589651982151578"
$ws.Range("E103").Value = [double]"0.1257413967680814"
$ws.Range("E105").Value = [double]"0.09097725694970631"
$ws.Range("E106").Value = [double]"2.514412390489136e-11"
$ws.Range("E107").Value = [double]"2.455480850087047e-14"
$ws.Range("E109").Value = [double]"0.04241100997540351"
$ws.Range("E110").Value = [double]"2.455480850087047e-14"
$ws.Range("E111").Value = [double]"0.09034800188302348"
$ws.Range("E112").Value = [double]"2.687769338505282e-10"
$ws.Range("E113").Value = [double]"0.06764313431491678"
$ws.Range("E114").Value = [double]"0.01291673091206783"
$ws.Range("E115").Value = [double]"0.1338724729757148"
$ws.Range("E116").Value = [double]"2.455480850087047e-14"
$ws.Range("E118").Value = [double]"0.003719439693224503"
$ws.Range("E119").Value = [double]"2.455480850087047e-14"
$ws.Range("E120").Value = [double]"2.455480850087047e-14"
$ws.Range("E121").Value = [double]"0.00371943969347005"
$ws.Range("E122").Value = [double]"2.455480850087047e-14"
$ws.Range("E123").Value = [double]"2.455480850087047e-14"
$ws.Range("E125").Value = [double]"2.455480850087047e-14"
$ws.Range("E127").Value = [double]"2.455480850087047e-14"
$ws.Range("E128").Value = [double]"2.455480850087047e-14"
$ws.Range("E129").Value = [double]"2.455480850087047e-14"
$ws.Range("E130").Value = [double]"0.01826999100572706"
$ws.Range("E131").Value = [double]"2.455480850087047e-14"
$ws.Range("E132").Value = [double]"2.455480850087047e-14"
$ws.Range("E133").Value = [double]"2.455480850087047e-14"
$ws.Range("E134").Value = [double]"2.514412390489136e-11"
$ws.Range("E135").Value = [double]"2.514412390489136e-11"
